$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = '21+0=21'
$t.Rows.Item(1).Cells.Item(2).Range.Text = '10+26=36'
$t.Rows.Item(1).Cells.Item(3).Range.Text = '88-65=23'
$t.Rows.Item(1).Cells.Item(4).Range.Text = '99-3=96'
$t.Rows.Item(1).Cells.Item(5).Range.Text = '68-56=12'
$t.Rows.Item(2).Cells.Item(1).Range.Text = '81-44=37'
$t.Rows.Item(2).Cells.Item(2).Range.Text = '5+57=62'
$t.Rows.Item(2).Cells.Item(3).Range.Text = '14+24=38'
$t.Rows.Item(2).Cells.Item(4).Range.Text = '58-1=57'
$t.Rows.Item(2).Cells.Item(5).Range.Text = '45+25=70'
$t.Rows.Item(3).Cells.Item(1).Range.Text = '92-52=40'
$t.Rows.Item(3).Cells.Item(2).Range.Text = '85-7=78'
$t.Rows.Item(3).Cells.Item(3).Range.Text = '47+24=71'
$t.Rows.Item(3).Cells.Item(4).Range.Text = '86-60=26'
$t.Rows.Item(3).Cells.Item(5).Range.Text = '35+62=97'
$t.Rows.Item(4).Cells.Item(1).Range.Text = '15+6=21'
$t.Rows.Item(4).Cells.Item(2).Range.Text = '13+52=65'
$t.Rows.Item(4).Cells.Item(3).Range.Text = '23-8=15'
$t.Rows.Item(4).Cells.Item(4).Range.Text = '47+25=72'
$t.Rows.Item(4).Cells.Item(5).Range.Text = '9+88=97'
$t.Rows.Item(5).Cells.Item(1).Range.Text = '69-29=40'
$t.Rows.Item(5).Cells.Item(2).Range.Text = '80+3=83'
$t.Rows.Item(5).Cells.Item(3).Range.Text = '65-30=35'
$t.Rows.Item(5).Cells.Item(4).Range.Text = '6+31=37'
$t.Rows.Item(5).Cells.Item(5).Range.Text = '8+20=28'
$t.Rows.Item(6).Cells.Item(1).Range.Text = '68-27=41'
$t.Rows.Item(6).Cells.Item(2).Range.Text = '25+52=77'
$t.Rows.Item(6).Cells.Item(3).Range.Text = '60-38=22'
$t.Rows.Item(6).Cells.Item(4).Range.Text = '58+18=76'
$t.Rows.Item(6).Cells.Item(5).Range.Text = '76-51=25'
$t.Rows.Item(7).Cells.Item(1).Range.Text = '48-35=13'
$t.Rows.Item(7).Cells.Item(2).Range.Text = '71-45=26'
$t.Rows.Item(7).Cells.Item(3).Range.Text = '53-22=31'
$t.Rows.Item(7).Cells.Item(4).Range.Text = '69-57=12'
$t.Rows.Item(7).Cells.Item(5).Range.Text = '2+86=88'
$t.Rows.Item(8).Cells.Item(1).Range.Text = '51-45=6'
$t.Rows.Item(8).Cells.Item(2).Range.Text = '11+45=56'
$t.Rows.Item(8).Cells.Item(3).Range.Text = '17+25=42'
$t.Rows.Item(8).Cells.Item(4).Range.Text = '23+9=32'
$t.Rows.Item(8).Cells.Item(5).Range.Text = '13+30=43'
$t.Rows.Item(9).Cells.Item(1).Range.Text = '65+32=97'
$t.Rows.Item(9).Cells.Item(2).Range.Text = '32+1=33'
$t.Rows.Item(9).Cells.Item(3).Range.Text = '33+11=44'
$t.Rows.Item(9).Cells.Item(4).Range.Text = '68-28=40'
$t.Rows.Item(9).Cells.Item(5).Range.Text = '29-3=26'
$t.Rows.Item(10).Cells.Item(1).Range.Text = '32+46=78'
$t.Rows.Item(10).Cells.Item(2).Range.Text = '10-6=4'
$t.Rows.Item(10).Cells.Item(3).Range.Text = '27+43=70'
$t.Rows.Item(10).Cells.Item(4).Range.Text = '1+35=36'
$t.Rows.Item(10).Cells.Item(5).Range.Text = '34+32=66'
$t.Rows.Item(11).Cells.Item(1).Range.Text = '99-9=90'
$t.Rows.Item(11).Cells.Item(2).Range.Text = '5+54=59'
$t.Rows.Item(11).Cells.Item(3).Range.Text = '95+4=99'
$t.Rows.Item(11).Cells.Item(4).Range.Text = '58-34=24'
$t.Rows.Item(11).Cells.Item(5).Range.Text = '7+59=66'
$t.Rows.Item(12).Cells.Item(1).Range.Text = '33+64=97'
$t.Rows.Item(12).Cells.Item(2).Range.Text = '98-57=41'
$t.Rows.Item(12).Cells.Item(3).Range.Text = '4+38=42'
$t.Rows.Item(12).Cells.Item(4).Range.Text = '1+64=65'
$t.Rows.Item(12).Cells.Item(5).Range.Text = '55-11=44'
$t.Rows.Item(13).Cells.Item(1).Range.Text = '16-13=3'
$t.Rows.Item(13).Cells.Item(2).Range.Text = '35+59=94'
$t.Rows.Item(13).Cells.Item(3).Range.Text = '20+17=37'
$t.Rows.Item(13).Cells.Item(4).Range.Text = '1+61=62'
$t.Rows.Item(13).Cells.Item(5).Range.Text = '10+48=58'
$t.Rows.Item(14).Cells.Item(1).Range.Text = '23+38=61'
$t.Rows.Item(14).Cells.Item(2).Range.Text = '89+5=94'
$t.Rows.Item(14).Cells.Item(3).Range.Text = '56-4=52'
$t.Rows.Item(14).Cells.Item(4).Range.Text = '19+47=66'
$t.Rows.Item(14).Cells.Item(5).Range.Text = '53-48=5'
$t.Rows.Item(15).Cells.Item(1).Range.Text = '77-74=3'
$t.Rows.Item(15).Cells.Item(2).Range.Text = '70-41=29'
$t.Rows.Item(15).Cells.Item(3).Range.Text = '37-19=18'
$t.Rows.Item(15).Cells.Item(4).Range.Text = '40+58=98'
$t.Rows.Item(15).Cells.Item(5).Range.Text = '90-45=45'
$t.Rows.Item(16).Cells.Item(1).Range.Text = '92-88=4'
$t.Rows.Item(16).Cells.Item(2).Range.Text = '42-40=2'
$t.Rows.Item(16).Cells.Item(3).Range.Text = '35+55=90'
$t.Rows.Item(16).Cells.Item(4).Range.Text = '17+1=18'
$t.Rows.Item(16).Cells.Item(5).Range.Text = '70-69=1'
$t.Rows.Item(17).Cells.Item(1).Range.Text = '79-50=29'
$t.Rows.Item(17).Cells.Item(2).Range.Text = '14+54=68'
$t.Rows.Item(17).Cells.Item(3).Range.Text = '92-13=79'
$t.Rows.Item(17).Cells.Item(4).Range.Text = '0+54=54'
$t.Rows.Item(17).Cells.Item(5).Range.Text = '40+55=95'
$t.Rows.Item(18).Cells.Item(1).Range.Text = '51-8=43'
$t.Rows.Item(18).Cells.Item(2).Range.Text = '96-21=75'
$t.Rows.Item(18).Cells.Item(3).Range.Text = '44+39=83'
$t.Rows.Item(18).Cells.Item(4).Range.Text = '57+1=58'
$t.Rows.Item(18).Cells.Item(5).Range.Text = '73+14=87'
$t.Rows.Item(19).Cells.Item(1).Range.Text = '30+59=89'
$t.Rows.Item(19).Cells.Item(2).Range.Text = '99-16=83'
$t.Rows.Item(19).Cells.Item(3).Range.Text = '87-1=86'
$t.Rows.Item(19).Cells.Item(4).Range.Text = '47+13=60'
$t.Rows.Item(19).Cells.Item(5).Range.Text = '9+54=63'
$t.Rows.Item(20).Cells.Item(1).Range.Text = '19+72=91'
$t.Rows.Item(20).Cells.Item(2).Range.Text = '51+5=56'
$t.Rows.Item(20).Cells.Item(3).Range.Text = '67-29=38'
$t.Rows.Item(20).Cells.Item(4).Range.Text = '40-33=7'
$t.Rows.Item(20).Cells.Item(5).Range.Text = '56-32=24'
